$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.071.44'
$ws.Range('E2').Value = '  +3.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.692.77'
$ws.Range('E3').Value = '  +2.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.52'
$ws.Range('E5').Value = '  +1.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.80'
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.580'
$ws.Range('E8').Value = '  +2.51%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.714.45'
$ws.Range('E9').Value = '  +2.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.64'
$ws.Range('E10').Value = '  +6.38%  '
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('E13').Value = '  +2.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.171.63'
$ws.Range('E14').Value = '  +2.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.964.66'
$ws.Range('E15').Value = '  +3.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.33'
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.702.42'
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '349.02'
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.57'
$ws.Range('E21').Value = '  +2.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.34'
$ws.Range('E22').Value = '  +2.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.82'
$ws.Range('E24').Value = '  +3.32%  '
$ws.Range('B25').Value = 'Polygon'
$ws.Range('C25').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.421'
$ws.Range('E25').Value = '  +0.37%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.171'
$ws.Range('E26').Value = '  +5.84%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.994'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0820'
$ws.Range('E28').Value = '  +1.99%  '
$ws.Range('E29').Value = '  +2.48%  '
$ws.Range('E30').Value = '  +9.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.61'
$ws.Range('E32').Value = '  +1.96%  '
$ws.Range('E33').Value = '  +0.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '150.02'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.26'
$ws.Range('E35').Value = '  +6.31%  '
$ws.Range('E36').Value = '  +8.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.948'
$ws.Range('E37').Value = '  -2.69%  '
$ws.Range('E38').Value = '  +4.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.53'
$ws.Range('E39').Value = '  +7.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.97'
$ws.Range('E40').Value = '  +0.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.67'
$ws.Range('E41').Value = '  -0.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '286.77'
$ws.Range('E42').Value = '  +3.09%  '
$ws.Range('E43').Value = '  +1.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.02'
$ws.Range('E44').Value = '  +2.04%  '
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.143.50'
$ws.Range('E47').Value = '  +7.75%  '
$ws.Range('E48').Value = '  +2.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0236'
$ws.Range('E49').Value = '  +2.69%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.81'
$ws.Range('E50').Value = '  +2.84%  '
$ws.Range('E51').Value = '  +1.67%  '
